$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.368.38'
$ws.Range("E2").Value = '  -1.18%  '

$ws.Range("D3").Value = '2.047.52'
$ws.Range("E3").Value = '  -1.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.86%  '

$ws.Range("E6").Value = '  -1.72%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.79'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.11%  '

$ws.Range("E9").Value = '  -2.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0815'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.84%  '

$ws.Range("E11").Value = '  -2.03%  '

$ws.Range("D12").Value = '2.348.62'
$ws.Range("E12").Value = '  -1.96%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.62'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.67'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.12%  '

$ws.Range("E15").Value = '  -3.17%  '

$ws.Range("E16").Value = '  -1.59%  '

$ws.Range("D17").Value = '2.049.06'
$ws.Range("E17").Value = '  -2.40%  '

$ws.Range("D18").Value = '37.200.04'
$ws.Range("E18").Value = '  -1.60%  '

$ws.Range("E19").Value = '  -1.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.86'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.96%  '

$ws.Range("D21").Value = '0.0₃0842'
$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.41'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("E24").Value = '  -1.08%  '

$ws.Range("E25").Value = '  -5.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.52'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.90'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.35%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.40'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.99'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.91%  '

$ws.Range("E31").Value = '  -2.77%  '

$ws.Range("E32").Value = '  -4.14%  '

$ws.Range("E33").Value = '  -1.97%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0613'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.40'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.83'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.995'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("E38").Value = '  -5.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("E40").Value = '  -7.02%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.481.83'
$ws.Range("E41").Value = '  +1.57%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.99'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.52%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.90'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0945'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.97'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.49%  '

$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("E47").Value = '  -4.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.11'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.94%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.235.89'
$ws.Range("E50").Value = '  -1.91%  '

$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.70'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -11.51%  '

